# BOT; UPDATE DATA
# Updates the kansensya_pcr.xlsx daily PCR-testing / infection tracking sheets:
#  - "all" sheet: last row (2020-05-11) test/positive/recovered counters bumped
#  - "kobe" sheet: new day of data rolled in (2020-05-12), tail rows revised
#  - the "market outside city" footnote is extended to include 2 more cases
#    (276 / 277) and its count 14/15 -> 16, and the two near-duplicate
#    footnotes collapse into a single shared string
#  - "kobe" becomes the active/visible tab again

$wb = $excel.ActiveWorkbook

$wsAll   = $wb.Worksheets.Item("all")
$wsKobe  = $wb.Worksheets.Item("kobe")
$wsOther = $wb.Worksheets.Item("other")

# ---------------------------------------------------------------------------
# 1) "all" sheet - row 34 (2020-05-11) revised totals
# ---------------------------------------------------------------------------
$wsAll.Cells.Item(34, 3).Value = 275   # C34 累計検査人数
$wsAll.Cells.Item(34, 4).Value = 81    # D34 陽性者数(累計)
$wsAll.Cells.Item(34, 5).Value = 70    # E34 治癒確認(累計)

# ---------------------------------------------------------------------------
# 2) "kobe" sheet - rows 79-89 revised cumulative test counts
# ---------------------------------------------------------------------------
$wsKobe.Cells.Item(79, 2).Value = 142
$wsKobe.Cells.Item(79, 3).Value = 2092

$wsKobe.Cells.Item(80, 3).Value = 2125
$wsKobe.Cells.Item(81, 3).Value = 2183
$wsKobe.Cells.Item(82, 3).Value = 2232
$wsKobe.Cells.Item(83, 3).Value = 2269
$wsKobe.Cells.Item(84, 3).Value = 2341
$wsKobe.Cells.Item(85, 3).Value = 2407
$wsKobe.Cells.Item(86, 3).Value = 2470
$wsKobe.Cells.Item(87, 3).Value = 2552
$wsKobe.Cells.Item(88, 3).Value = 2590

# New day rolled in: 2020-05-12
$wsKobe.Cells.Item(89, 2).Value = 69    # B89 本日検査人数
$wsKobe.Cells.Item(89, 3).Value = 2659  # C89 累計検査人数
$wsKobe.Cells.Item(89, 6).Value = 76    # F89
$wsKobe.Cells.Item(89, 7).Value = 66    # G89

# ---------------------------------------------------------------------------
# 3) Footnote about out-of-city residents: extend the case list and merge
#    the two near-duplicate shared strings into a single, reused one.
# ---------------------------------------------------------------------------
$footnote = "※　24・34・53・58・59・60・158・161・163・192・237・248・268・272・276・277例目（計16件）は市外在住者です。"
$wsAll.Cells.Item(35, 2).Value = $footnote
$wsKobe.Cells.Item(90, 2).Value = $footnote

# ---------------------------------------------------------------------------
# 4) View state: "kobe" becomes the active/visible sheet and tab again,
#    selections are moved to match, and the kobe sheet zoom goes 70% -> 85%.
# ---------------------------------------------------------------------------
$wsAll.Range("A33").Select()
$wsOther.Range("J64").Select()
$wsKobe.Range("B89").Select()
$excel.ActiveWindow.Zoom = 85
